$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-07-04 Friday" "2025-07-05 Saturday"

Replace-Text "226÷4=56, 2" "781÷8=97, 5"
Replace-Text "332÷4=83, 0" "541÷4=135, 1"
Replace-Text "353÷6=58, 5" "774÷9=86, 0"
Replace-Text "877÷8=109, 5" "987÷8=123, 3"
Replace-Text "225÷5=45, 0" "542÷2=271, 0"

Replace-Text "647÷2=323, 1" "291÷9=32, 3"
Replace-Text "230÷2=115, 0" "488÷6=81, 2"
Replace-Text "960÷2=480, 0" "474÷9=52, 6"
Replace-Text "550÷7=78, 4" "934÷9=103, 7"
Replace-Text "894÷2=447, 0" "604÷7=86, 2"

Replace-Text "374÷4=93, 2" "158÷3=52, 2"
Replace-Text "567÷7=81, 0" "850÷4=212, 2"
Replace-Text "759÷7=108, 3" "211÷9=23, 4"
Replace-Text "734÷6=122, 2" "576÷4=144, 0"
Replace-Text "119÷8=14, 7" "945÷8=118, 1"

Replace-Text "436÷8=54, 4" "842÷9=93, 5"
Replace-Text "619÷3=206, 1" "462÷4=115, 2"
Replace-Text "536÷6=89, 2" "267÷3=89, 0"
Replace-Text "874÷2=437, 0" "353÷7=50, 3"
Replace-Text "620÷2=310, 0" "486÷9=54, 0"

Replace-Text "536÷2=268, 0" "713÷4=178, 1"
Replace-Text "502÷6=83, 4" "465÷5=93, 0"
Replace-Text "256÷2=128, 0" "882÷2=441, 0"
Replace-Text "202÷9=22, 4" "663÷5=132, 3"
Replace-Text "573÷6=95, 3" "614÷3=204, 2"
